# get highest row and col -> read all
#
# Probe the sheet for its highest used row/column (mirrors a reader that
# first finds the extent of the table, then walks every cell in it), then
# extend that table: column F switches from a computed floating point
# number to a pre-formatted decimal string ("1,23" style), and five more
# data rows are appended below the existing ones.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used     = $ws.UsedRange
$firstRow = $used.Row
$firstCol = $used.Column
$lastRow  = $firstRow + $used.Rows.Count - 1
$lastCol  = $firstCol + $used.Columns.Count - 1

$headerRow = $firstRow
$dataFirstRow = $headerRow + 1
$dataLastRow  = $lastRow

$colA = $firstCol
$colF = $colA + 5

$extraRows = 5
$newLastRow = $dataLastRow + $extraRows

# "1,23" .. "1,33" - one value per data row, old and new alike.
$n = $dataLastRow - $dataFirstRow + 1 + $extraRows
$fText = New-Object 'object[]' $n
for ($i = 0; $i -lt $n; $i++) {
    $hundredths = 23 + $i
    $fText[$i] = "1," + $hundredths
}

# Labels for the newly appended rows only.
$newLabels = @("string 7", "string 8", "string 9", "string 10", "string 11")

for ($r = $dataFirstRow; $r -le $newLastRow; $r++) {
    $idx = $r - $dataFirstRow

    if ($r -gt $dataLastRow) {
        $ws.Cells.Item($r, $colA).Value = $newLabels[$r - $dataLastRow - 1]
    }

    $fCell = $ws.Cells.Item($r, $colF)
    $fCell.NumberFormat = "0.00"
    $fCell.Value = $fText[$idx]
}
